{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 ... Creative Commons Attribution\" line, and the blank\n// paragraph that separates them from the \"LOB1037: ...\" requirement line\n// above (the page footer/boilerplate block Jekyll injects at build time).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nlet reqIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"LOB1037\") !== -1) {\n    reqIndex = i;\n  }\n  if (text.indexOf(\"Creative Commons Attribution\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (reqIndex !== -1 && copyrightIndex !== -1 && copyrightIndex > reqIndex) {\n  // Delete from the end backward so earlier indices stay valid.\n  for (let i = copyrightIndex; i > reqIndex; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the three paragraphs that are being removed from the page\n# footer block: the blank paragraph right after the \"LOB1037: ...\"\n# requirement line, the \"Ver no Jupiter ...\" line, and the\n# \"(c) 2020 ... Creative Commons Attribution\" line.\n$startPara = $null\n$endPara = $null\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text\n\n    if ($text -match \"LOB1037\") {\n        # The paragraph immediately following the LOB1037 requirement\n        # line is the blank paragraph that starts the block to delete.\n        $startPara = $d.Paragraphs.Item($i + 1)\n    }\n\n    if ($text -match [regex]::Escape(\"Creative Commons Attribution\")) {\n        $endPara = $para\n    }\n}\n\nif ($startPara -ne $null -and $endPara -ne $null) {\n    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $deleteRange.Delete()\n}\n"}
